# Re-ran "resolve" and "classify+summarise" steps after changes to mapping file.
# This updates the summary sheets: several "Range Status" and "Species
# qualification" figures now compute to 0 (no species left mapped to a
# state/range after the mapping-file change), and the "High Priority
# break-up" sheet now only has two break-up categories (Trend New / IUCN)
# instead of four, with their counts/percentages recomputed.

$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": B column -> 0, remove C column (perc.) entirely ---
$ws2 = $wb.Worksheets.Item("Range Status")
foreach ($r in 2..7) {
    $ws2.Cells.Item($r, 2).Value = 0
    $ws2.Cells.Item($r, 3).ClearContents()
}

# --- Sheet "Species qualification": Range Analysis species count -> 0 ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B5").Value = 0

# --- Sheet "High Priority break-up": recompute rows, drop Range & old IUCN rows ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 3
$ws5.Range("C2").Value = 15.8
$ws5.Range("D2").Value = 3
$ws5.Range("E2").Value = 15.8

$ws5.Range("A3").Value = "IUCN"
$ws5.Range("B3").Value = 16
$ws5.Range("C3").Value = 84.2
$ws5.Range("D3").Value = 16
$ws5.Range("E3").Value = 84.2

# remove the now-obsolete "Range" and old "IUCN" rows
$ws5.Rows("4:5").Delete()
